# Hypothesis testing example added
# Fix typo in the "gender" row comment: "Makes" -> "Males"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "1. Two distinct values`n2. Males are more than females`n3. 23% missing values"

# Update the active view/selection to match the post-edit state
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("B13").Select()
